# Auto update Excel log
# Appends newly-ingested sensor events to the ALERTS log and the mmWave
# (Living Room) log.
#
# Date-looking strings (column A, e.g. "2026-02-01") get forced to Text
# number format first so Excel's COM layer stores them as literal strings
# instead of auto-converting them to date serials.

$wb = $excel.ActiveWorkbook

# --- ALERTS sheet: append row 2 (a new CRITICAL fall-detection alert) ---
$alerts = $wb.Worksheets.Item("ALERTS")

$alerts.Range("A2").NumberFormat = "@"
$alerts.Range("A2").Value = "2026-02-01"
$alerts.Range("B2").Value = "17:21:22"
$alerts.Range("C2").Value = "17:00"
$alerts.Range("D2").Value = "Living Room"
$alerts.Range("E2").Value = "CRITICAL"
$alerts.Range("F2").Value = "FALL_DETECTED"

# --- mmWave sheet: append rows 31-33 (presence detections) ---
$mmwave = $wb.Worksheets.Item("mmWave")

$mmwaveRows = @(
    @{ Row = 31; Time = "17:21:03" },
    @{ Row = 32; Time = "17:21:47" },
    @{ Row = 33; Time = "17:21:57" }
)

foreach ($entry in $mmwaveRows) {
    $r = $entry.Row
    $mmwave.Range("A$r").NumberFormat = "@"
    $mmwave.Range("A$r").Value = "2026-02-01"
    $mmwave.Range("B$r").Value = $entry.Time
    $mmwave.Range("C$r").Value = "17:00"
    $mmwave.Range("D$r").Value = "Living Room"
    $mmwave.Range("E$r").Value = "PRESENCE_DETECTED"
    $mmwave.Range("F$r").Value = "Active"
}
